# Implement 5-cell professional notebook structure
# 1. Insert a new "Time Analysis" worksheet after "Profitability" (and before
#    "Branch Payment Pref"), populated with time-of-day sales data.
# 2. Register the new analysis as a row in the "Executive Summary" sheet,
#    inserted before the "Branch Payment Method Preferences" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the "Time Analysis" worksheet in the correct tab position
# ---------------------------------------------------------------------
$profitSheet = $wb.Worksheets.Item("Profitability")
$timeSheet = $wb.Worksheets.Add($null, $profitSheet)
$timeSheet.Name = "Time Analysis"

# Header row
$timeSheet.Range("A1").Value = "time_period"
$timeSheet.Range("B1").Value = "total_transactions"
$timeSheet.Range("C1").Value = "total_revenue"
$timeSheet.Range("D1").Value = "avg_transaction_value"

# Match the bold / centered / thin-bordered header style used on the other
# sheets in the workbook.
$headerRange = $timeSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$timeSheet.Range("A2").Value = "Afternoon"
$timeSheet.Range("B2").Value = 4636
$timeSheet.Range("C2").Value = 578817.77
$timeSheet.Range("D2").Value = 124.85

$timeSheet.Range("A3").Value = "Evening"
$timeSheet.Range("B3").Value = 3246
$timeSheet.Range("C3").Value = 379401.6
$timeSheet.Range("D3").Value = 116.88

$timeSheet.Range("A4").Value = "Morning"
$timeSheet.Range("B4").Value = 2087
$timeSheet.Range("C4").Value = 251507.01
$timeSheet.Range("D4").Value = 120.51

# ---------------------------------------------------------------------
# 2. Insert the matching catalog row into "Executive Summary"
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("Executive Summary")
$summarySheet.Range("A6:D6").Insert()

$summarySheet.Range("A6").Value = "Sales Performance by Time of Day"
$summarySheet.Range("B6").Value = 3
$summarySheet.Range("C6").Value = "Operational efficiency and staff scheduling optimization"
$summarySheet.Range("D6").Value = "6_sales_by_time_period.csv"

Write-Output "Time Analysis sheet inserted and Executive Summary updated"
